$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Difference to predicted" header cell in H26 ---------------------
# Pull border (left+right only, matches H10/H15/H16/H19 "boxed value" style)
$ws.Range("H10").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
# Pull the blue fill used by the row 26/27 header band
$ws.Range("A24").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
# Re-apply the thin left/right border (the fill paste above reset it)
$ws.Range("H26").Borders.Item(7).LineStyle = 1
$ws.Range("H26").Borders.Item(7).Weight = 2
$ws.Range("H26").Borders.Item(10).LineStyle = 1
$ws.Range("H26").Borders.Item(10).Weight = 2

$ws.Range("H26").Value = "Difference to predicted"
$ws.Range("H26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 45

# --- New "Difference to predicted" values for the MANUAL COUNT table ------
$ws.Range("H28").Value = 0.2
$ws.Range("H29").Value = 0.9

$ws.Range("E30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = 0.5

$ws.Range("E31").Copy() | Out-Null
$ws.Range("H31").PasteSpecial(-4122) | Out-Null
$ws.Range("H31").Value = 1.2

$ws.Range("E32").Copy() | Out-Null
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("H32").Value = 0

$ws.Range("E33").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null
$ws.Range("H33").Value = 0.1

# --- Summary row 34: per-column averages -----------------------------------
$ws.Range("B34").Formula = "=(27+19+18+12+10+8)/6"
$ws.Range("C34").Formula = "=(31+29+21+19+14+11)/6"
$ws.Range("D34").Formula = "=(4+9+14+17+23+28)/6"
$ws.Range("E34").Formula = "=(3+7+10+12+15+21)/6"
$ws.Range("H34").Formula = "=AVERAGE(H28:H33)"

# --- New row 35: combined OFF/ON averages -----------------------------------
$ws.Range("B35").Formula = "=AVERAGE(B34:C34)"
$ws.Range("D35").Formula = "=AVERAGE(D34:E34)"

# --- Selection matches the new last-used cell -------------------------------
$ws.Range("H35").Select() | Out-Null
